$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 1739.375
$ws.Range("I11").Value = 1739.375
$ws.Range("K11").Value = 1739.375
$ws.Range("M11").Value = -1599.375

# Row 51
$ws.Range("H51").Value = 14245.3125
$ws.Range("I51").Value = 20297.4
$ws.Range("J51").Value = 13124.556
$ws.Range("K51").Value = 20297.4
$ws.Range("L51").Value = 13124.556
$ws.Range("M51").Value = -19813.4
$ws.Range("N51").Value = -14092.556

# Row 53
$ws.Range("H53").Value = 930
$ws.Range("I53").Value = 736.8570999999999
$ws.Range("K53").Value = 736.8570999999999
$ws.Range("M53").Value = -99.85709999999995

# Row 88
$ws.Range("H88").Value = 3399.9092
$ws.Range("I88").Value = 3216.6667
$ws.Range("J88").Value = 3619.8
$ws.Range("K88").Value = 3216.6667
$ws.Range("L88").Value = 3619.8
$ws.Range("M88").Value = -2810.6667
$ws.Range("N88").Value = -4431.8

# Row 91
$ws.Range("H91").Value = 3399.9092
$ws.Range("I91").Value = 3216.6667
$ws.Range("J91").Value = 3619.8
$ws.Range("K91").Value = 3216.6667
$ws.Range("L91").Value = 3619.8
$ws.Range("M91").Value = -1812.6667
$ws.Range("N91").Value = -6427.8

# Row 92
$ws.Range("H92").Value = 2629.682
$ws.Range("I92").Value = 2521.5334
$ws.Range("K92").Value = 2521.5334
$ws.Range("M92").Value = -1273.5334

# Row 111
$ws.Range("H111").Value = 1001.4545
$ws.Range("I111").Value = 968.7778
$ws.Range("K111").Value = 2906.3334
$ws.Range("M111").Value = 160.6666

# Row 132
$ws.Range("H132").Value = 7281.769
$ws.Range("I132").Value = 6687.125
$ws.Range("K132").Value = 20061.375
$ws.Range("M132").Value = -17531.375

# Row 138
$ws.Range("H138").Value = 2692.875
$ws.Range("I138").Value = 1251.3478
$ws.Range("J138").Value = 3202.9539
$ws.Range("K138").Value = 3754.0434
$ws.Range("L138").Value = 9608.861699999999
$ws.Range("M138").Value = 1385.9566
$ws.Range("N138").Value = -19888.8617

# Row 139
$ws.Range("H139").Value = 99933
$ws.Range("J139").Value = 99933
$ws.Range("L139").Value = 99933
$ws.Range("N139").Value = -110213

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10398.342
$ws.Range("I32").Value = 8012.543
$ws.Range("J32").Value = 38232.668
$ws.Range("K32").Value = 8012.543
$ws.Range("L32").Value = 38232.668
$ws.Range("M32").Value = -7725.543
$ws.Range("N32").Value = -38806.668

# Row 122
$ws.Range("H122").Value = 2887.9333
$ws.Range("I122").Value = 1981.2222
$ws.Range("J122").Value = 4248
$ws.Range("K122").Value = 5943.6666
$ws.Range("L122").Value = 12744
$ws.Range("M122").Value = -3493.6666
$ws.Range("N122").Value = -17644

# Row 135
$ws.Range("H135").Value = 64142.332
$ws.Range("J135").Value = 64142.332
$ws.Range("L135").Value = 64142.332
$ws.Range("N135").Value = -74282.33199999999

# Row 139
$ws.Range("H139").Value = 99999
$ws.Range("J139").Value = 99999
$ws.Range("L139").Value = 99999
$ws.Range("N139").Value = -110279

$ws = $wb.Worksheets.Item("BSM")
# Row 81
$ws.Range("H81").Value = 50339
$ws.Range("J81").Value = 50339
$ws.Range("L81").Value = 50339
$ws.Range("N81").Value = -52461

# Row 84
$ws.Range("H84").Value = 50339
$ws.Range("J84").Value = 50339
$ws.Range("L84").Value = 151017
$ws.Range("N84").Value = -161625

# Row 99
$ws.Range("H99").Value = 5012.75
$ws.Range("I99").Value = 5129
$ws.Range("J99").Value = 4757
$ws.Range("K99").Value = 5129
$ws.Range("L99").Value = 4757
$ws.Range("M99").Value = -3631
$ws.Range("N99").Value = -7753

$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 3499.6667
$ws.Range("I10").Value = 3499.6667
$ws.Range("K10").Value = 3499.6667
$ws.Range("M10").Value = -3360.6667

# Row 31
$ws.Range("H31").Value = 359040.72
$ws.Range("I31").Value = 589198.1
$ws.Range("J31").Value = 3342.9092
$ws.Range("K31").Value = 589198.1
$ws.Range("L31").Value = 3342.9092
$ws.Range("M31").Value = -588903.1
$ws.Range("N31").Value = -3932.9092

# Row 34
$ws.Range("H34").Value = 359040.72
$ws.Range("I34").Value = 589198.1
$ws.Range("J34").Value = 3342.9092
$ws.Range("K34").Value = 589198.1
$ws.Range("L34").Value = 3342.9092
$ws.Range("M34").Value = -588996.1
$ws.Range("N34").Value = -3746.9092

# Row 64
$ws.Range("H64").Value = 35838.2
$ws.Range("J64").Value = 35838.2
$ws.Range("L64").Value = 35838.2
$ws.Range("N64").Value = -36334.2

# Row 67
$ws.Range("H67").Value = 35838.2
$ws.Range("J67").Value = 35838.2
$ws.Range("L67").Value = 35838.2
$ws.Range("N67").Value = -37554.2

# Row 122
$ws.Range("H122").Value = 2337.0908
$ws.Range("I122").Value = 1960.2
$ws.Range("J122").Value = 2651.1667
$ws.Range("K122").Value = 5880.6
$ws.Range("L122").Value = 7953.500100000001
$ws.Range("M122").Value = -3430.6
$ws.Range("N122").Value = -12853.5001

# Row 132
$ws.Range("H132").Value = 2297.7273
$ws.Range("I132").Value = 2320.4443
$ws.Range("K132").Value = 6961.3329
$ws.Range("M132").Value = -4431.3329

# Row 134
$ws.Range("H134").Value = 7757.125
$ws.Range("I134").Value = 8714.895
$ws.Range("K134").Value = 26144.685
$ws.Range("M134").Value = -23609.685

$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 1999.5
$ws.Range("I22").Value = 2000
$ws.Range("K22").Value = 6000
$ws.Range("M22").Value = -5831

# Row 27
$ws.Range("H27").Value = 1999.5
$ws.Range("I27").Value = 2000
$ws.Range("K27").Value = 6000
$ws.Range("M27").Value = -5898

# Row 40
$ws.Range("H40").Value = 1002
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1002
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 4008
$ws.Range("M40").Value = ""
$ws.Range("N40").Value = -4146

# Row 50
$ws.Range("H50").Value = 949.5
$ws.Range("I50").Value = 899
$ws.Range("K50").Value = 2697
$ws.Range("M50").Value = -2216

# Row 53
$ws.Range("H53").Value = 949.5
$ws.Range("I53").Value = 899
$ws.Range("K53").Value = 2697
$ws.Range("M53").Value = -2216

# Row 68
$ws.Range("H68").Value = 13890422
$ws.Range("I68").Value = 41667268
$ws.Range("J68").Value = 1999.75
$ws.Range("K68").Value = 125001804
$ws.Range("L68").Value = 5999.25
$ws.Range("M68").Value = -125000993
$ws.Range("N68").Value = -7621.25

# Row 71
$ws.Range("H71").Value = 13890422
$ws.Range("I71").Value = 41667268
$ws.Range("J71").Value = 1999.75
$ws.Range("K71").Value = 375005412
$ws.Range("L71").Value = 17997.75
$ws.Range("M71").Value = -375001356
$ws.Range("N71").Value = -26109.75

# Row 75
$ws.Range("H75").Value = 608.6667
$ws.Range("I75").Value = 608.6667
$ws.Range("K75").Value = 1826.0001
$ws.Range("M75").Value = -828.0001

# Row 78
$ws.Range("H78").Value = 608.6667
$ws.Range("I78").Value = 608.6667
$ws.Range("K78").Value = 5478.0003
$ws.Range("M78").Value = -486.0002999999997

# Row 98
$ws.Range("H98").Value = 387.69232
$ws.Range("I98").Value = 402
$ws.Range("J98").Value = 375.42856
$ws.Range("K98").Value = 1206
$ws.Range("L98").Value = 1126.28568
$ws.Range("M98").Value = 292
$ws.Range("N98").Value = -4122.28568

# Row 103
$ws.Range("H103").Value = 1561.8572
$ws.Range("I103").Value = 911.3333
$ws.Range("J103").Value = 2049.75
$ws.Range("K103").Value = 2733.9999
$ws.Range("L103").Value = 6149.25
$ws.Range("M103").Value = -1854.9999
$ws.Range("N103").Value = -7907.25

# Row 114
$ws.Range("H114").Value = 605.625
$ws.Range("J114").Value = 682.5
$ws.Range("L114").Value = 2047.5
$ws.Range("N114").Value = -8555.5

# Row 117
$ws.Range("H117").Value = 483.4
$ws.Range("I117").Value = 299.33334
$ws.Range("J117").Value = 759.5
$ws.Range("K117").Value = 898.0000200000001
$ws.Range("L117").Value = 2278.5
$ws.Range("M117").Value = 2543.99998
$ws.Range("N117").Value = -9162.5

# Row 134
$ws.Range("H134").Value = 8958.786
$ws.Range("I134").Value = 11183.692
$ws.Range("K134").Value = 33551.076
$ws.Range("M134").Value = -28481.076

# Row 138
$ws.Range("H138").Value = 1964.9
$ws.Range("I138").Value = 1964.9
$ws.Range("K138").Value = 5894.700000000001
$ws.Range("M138").Value = -754.7000000000007

# Row 139
$ws.Range("H139").Value = 3109.4783
$ws.Range("I139").Value = 2352.0588
$ws.Range("K139").Value = 7056.176399999999
$ws.Range("M139").Value = -1916.176399999999

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 7448.1313
$ws.Range("J70").Value = 7433.1763
$ws.Range("L70").Value = 7433.1763
$ws.Range("N70").Value = -7973.1763

# Row 73
$ws.Range("H73").Value = 7448.1313
$ws.Range("J73").Value = 7433.1763
$ws.Range("L73").Value = 7433.1763
$ws.Range("N73").Value = -9305.176299999999

# Row 102
$ws.Range("H102").Value = 3324.2917
$ws.Range("I102").Value = 3286.5454
$ws.Range("J102").Value = 3739.5
$ws.Range("K102").Value = 3286.5454
$ws.Range("L102").Value = 3739.5
$ws.Range("M102").Value = -1664.5454
$ws.Range("N102").Value = -6983.5

# Row 122
$ws.Range("H122").Value = 4576.222
$ws.Range("I122").Value = 9130.666999999999
$ws.Range("K122").Value = 27392.001
$ws.Range("M122").Value = -24942.001

# Row 126
$ws.Range("H126").Value = 10472.75
$ws.Range("J126").Value = 4064.75
$ws.Range("L126").Value = 12194.25
$ws.Range("N126").Value = -17134.25

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1475.7142
$ws.Range("J16").Value = 4048.5
$ws.Range("L16").Value = 4048.5
$ws.Range("N16").Value = -4388.5

$ws = $wb.Worksheets.Item("WVR")
# Row 47
$ws.Range("H47").Value = 39503.89
$ws.Range("J47").Value = 31070
$ws.Range("L47").Value = 31070
$ws.Range("N47").Value = -32214

# Row 122
$ws.Range("H122").Value = 19679.5
$ws.Range("I122").Value = 2599.077
$ws.Range("K122").Value = 7797.231000000001
$ws.Range("M122").Value = -5347.231000000001

# Row 126
$ws.Range("H126").Value = 3253.8333
$ws.Range("I126").Value = 3204.0588
$ws.Range("K126").Value = 9612.1764
$ws.Range("M126").Value = -7142.1764
